# "alteração no arquivo de ajuda"
#
# 1. Insert a new first paragraph at the very top of the document:
#      "TESTE EFETUADO USANDO O AMBIENTE WINDOWS"
#    formatted with a 1416-twip (70.8pt) left indent and 14pt (sz=28 half-pts)
#    text, using the plain "Normal" style (no list numbering).
#
# 2. Drop the stray <w:lastRenderedPageBreak/> marker that sits in front of
#    the "Rodar a aplicação PRODUTORA..." bullet - delete + retype that run's
#    text so the render-break bookmark isn't carried over.
#
# 3. (styles.xml) A <w:rsid w:val="00EC78A9"/> is stamped onto the "Normal"
#    style definition by Word's own RSID bookkeeping for this editing
#    session; the COM surface exposed here has no property for it, so it is
#    not reproducible through automation and is intentionally left alone.

$d = $word.ActiveDocument

# --- 1. Insert the new heading paragraph at the very top of the document ---
$firstPar = $d.Paragraphs.Item(1)
$firstPar.Range.InsertParagraphBefore()

$newPar = $d.Paragraphs.Item(1)
$newPar.Style = "Normal"
$newPar.Range.ParagraphFormat.LeftIndent = 70.8
$newPar.Range.Text = "TESTE EFETUADO USANDO O AMBIENTE WINDOWS"
$newPar.Range.Font.Size = 14

# --- 2. Remove the lastRenderedPageBreak marker on the PRODUTORA bullet ---
$produtoraPar = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*PRODUTORA*") {
        $produtoraPar = $p
        break
    }
}

if ($produtoraPar -ne $null) {
    $r = $produtoraPar.Range
    [void]$r.MoveEnd(1, -1)          # exclude the trailing paragraph mark
    $text = $r.Text
    [void]$r.Delete()
    [void]$r.InsertAfter($text)
}
